# Apply the "PO Forecast" update:
#  1. Rename the "Requested quantity" headers on the existing sheets.
#  2. Add a new "PO Forecast" worksheet at the end of the workbook.
#  3. Populate it with the forecast header row + 19 data rows, matching
#     the formatting used on the existing sheets (bold/bordered header,
#     date-formatted first column).

$wb = $excel.ActiveWorkbook

# --- 1. Rename the quantity headers -----------------------------------
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet --
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Header row ------------------------------------------------------
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header formatting (bold, centered, thin border) used elsewhere.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# --- 4. Data rows ---------------------------------------------------
$rows = @(
    @(2, 45088.99999999999, 4, -53.59952057100707, 63.35373748978834),
    @(3, 45095.99999999999, 6, -50.03684270691731, 64.71401748175397),
    @(4, 45102.99999999999, 9, -50.26151999735841, 63.06010831526781),
    @(5, 45116.99999999999, 13, -45.36418518728878, 71.55786073631728),
    @(6, 45186.99999999999, 35, -20.14614481566314, 92.94933173830368),
    @(7, 45193.99999999999, 37, -20.38042305636975, 98.88023779195989),
    @(8, 45207.99999999999, 42, -20.00053648462732, 95.42021331264914),
    @(9, 45214.99999999999, 44, -18.24402246857701, 97.95160560944097),
    @(10, 45221.99999999999, 46, -16.14065100065714, 99.87610065152549),
    @(11, 45242.99999999999, 53, -5.586819940897036, 107.5992524410252),
    @(12, 45249.99999999999, 55, -4.254887398829669, 115.8515952388666),
    @(13, 45256.99999999999, 57, -1.334890218979557, 110.9869417929817),
    @(14, 45263.99999999999, 59, 6.143501711016082, 117.2799618864204),
    @(15, 45270.99999999999, 61, 5.699137713993807, 116.3632502067352),
    @(16, 45277.99999999999, 64, 4.104072447334478, 124.6142034848946),
    @(17, 45284.99999999999, 66, 9.986881054697804, 125.6500990927476),
    @(18, 45291.99999999999, 68, 2.70791775408112, 120.9206100984725),
    @(19, 45298.99999999999, 70, 10.35245978567983, 128.0372536278802),
    @(20, 45305.99999999999, 72, 13.10252086951075, 128.0804699188382)
)

foreach ($row in $rows) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value = $row[1]
    $wsForecast.Cells.Item($r, 2).Value = $row[2]
    $wsForecast.Cells.Item($r, 3).Value = $row[3]
    $wsForecast.Cells.Item($r, 4).Value = $row[4]
}

# Apply the date number format (and matching style) used for column A on
# the other sheets to the new "ds" column's data cells.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A20").PasteSpecial(-4122)

# Re-set the data row A values after the paste (PasteSpecial with the
# formats-only option keeps the numeric values already entered, but make
# sure nothing was clobbered).
foreach ($row in $rows) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value = $row[1]
}
